$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("H3").Value = 0.7935955447267664
$ws.Range("I3").Value = 0.06059009483667018
$ws.Range("K3").Value = 186
$ws.Range("Q3").Value = 28
$ws.Range("S3").Value = 83
$ws.Range("T3").Value = 127
$ws.Range("U3").Value = 239
$ws.Range("V3").Value = 2819
$ws.Range("X3").Value = 2764
$ws.Range("Y3").Value = 2720
$ws.Range("Z3").Value = 2608
$ws.Range("AF3").Value = 0.990165
$ws.Range("AH3").Value = 0.970847
$ws.Range("AI3").Value = 0.955392
$ws.Range("AJ3").Value = 0.916052
